$d = $word.ActiveDocument

# Locate the start of the paragraph's text ("Esteu participant...")
$startRng = $d.Content
$startRng.Find.Execute("Esteu participant")
$startPos = $startRng.Start

# Locate the end of the paragraph's text (the final sentence)
$endRng = $d.Content
$endRng.Find.Execute("faran palesa la visibilitat del cel nocturn.")
$endPos = $endRng.End

# Replace the whole run sequence (which spans many small w:r elements)
# with a single merged run containing the updated wording.
$fullRange = $d.Range($startPos, $endPos)
$fullRange.Delete()

$insRange = $d.Range($startPos, $startPos)
$insRange.InsertAfter("Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles més febles que es poden veure, com a mitjà per mesurar la contaminació lumínica en un lloc determinat. Localitzant i observant la  constel·lació d'Orió a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el món aprendran com els llums de la seva zona contribueixen a augmentar la contaminació lumínica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn.")
